$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datos")

# Row 10: numeric hours worked for Semanas 9-12
$ws.Range("B10").Value = 11
$ws.Range("C10").Value = 9.5
$ws.Range("D10").Value = 21
$ws.Range("E10").Value = 12
$ws.Range("F10").Value = 5
$ws.Range("G10").Value = 5
$ws.Range("H10").Value = 4
$ws.Range("I10").Value = 5.5
$ws.Range("J10").Value = 11.5

# Row 11: work descriptions for Semanas 9-12
# NOTE: assignment order matters — it controls the order new entries land in
# xl/sharedStrings.xml (and therefore their <si> index), so cells are written
# in the same first-seen order as the authoring session (B,J,E,H,D,F,G,C,I).
$ws.Range("B11").Value = "1h maquetación, `n7h documentación`n1h reunión`n1h retoque web`n1h repaso de estado y org."
$ws.Range("J11").Value = "1,5 reuniones`n20min actas`n40min documentación`n2,5h revisión guia instalación`n2 h implementación BBDD`n4,5h población final BBDD"
$ws.Range("E11").Value = "4h Implementación`n1h poblado de BBDD`n6h poblacion final BBDD`n1h pruebas"
$ws.Range("H11").Value = "2,5 reuniones`n1,5 documentacion"
$ws.Range("D11").Value = "2h reunion`n3h funcionalidades de análisis`n1h comentarios`n2h bug fix y pruebas`n1h temas`n1,5h calificaciones`n1,5h modulo`n1h routers perfil privado`n1h eliminar perfil"
$ws.Range("F11").Value = "1h reunion`n3h documentacion"
$ws.Range("G11").Value = "30 min organización reuniones`n3h documentacion`n1h Reunión`n0,5h manual de usuario"
$ws.Range("C11").Value = "2h reunion`n1h documentacion`n1,5 Scripts de instalacion`n2h Perfil publico`n3h finalizar perfil publico y privado"
$ws.Range("I11").Value = "1,5 reuniones`n20min actas`n2,66 documentación"

# Update the view: scroll to A13, select I13
$ws.Application.ActiveWindow.ScrollRow = 13
$ws.Range("I13").Select()
